# Fruta / hortaliza, semanal
# Insert one new week of price records (5 rows) for "Zapallo italiano" at
# Vega Modelo de Temuco, right before the existing row 233, pushing the old
# rows 233-292 down to 238-297 (dimension grows from A1:R292 to A1:R297).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 233:292 down by 5 rows, creating 5 blank rows at 233:237.
$ws.Rows("233:237").Insert()

# Data for the new rows (common columns across this sub-table).
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 100112032
$categoria = "Zapallo italiano"
$clasif    = "Hortaliza"

function Set-Fila {
    param(
        [int]$Fila,
        [double]$Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidades
    )

    $ws.Cells.Item($Fila, 1).Value = $mercadoId
    $ws.Cells.Item($Fila, 2).Value = $mercado
    $ws.Cells.Item($Fila, 3).Value = $region
    $ws.Cells.Item($Fila, 4).Value = $Fecha
    $ws.Cells.Item($Fila, 5).Value = $codreg
    $ws.Cells.Item($Fila, 6).Value = $catId
    $ws.Cells.Item($Fila, 7).Value = $categoria
    $ws.Cells.Item($Fila, 8).Value = $Variedad
    $ws.Cells.Item($Fila, 9).Value = $Calidad
    $ws.Cells.Item($Fila, 10).Value = $Volumen
    $ws.Cells.Item($Fila, 11).Value = $PrecioMin
    $ws.Cells.Item($Fila, 12).Value = $PrecioMax
    $ws.Cells.Item($Fila, 13).Value = $PrecioProm
    $ws.Cells.Item($Fila, 14).Value = $Unidad
    $ws.Cells.Item($Fila, 15).Value = $Origen
    $ws.Cells.Item($Fila, 16).Value = $PrecioKg
    $ws.Cells.Item($Fila, 17).Value = $KgUnidades
    $ws.Cells.Item($Fila, 18).Value = $clasif
}

Set-Fila -Fila 233 -Fecha 44511 -Variedad "Bola 8"            -Calidad "Primera" -Volumen 80  -PrecioMin 12000 -PrecioMax 12000 -PrecioProm 12000 -Unidad "`$/caja 60 unidades" -Origen "Región de O'Higgins"          -PrecioKg 200 -KgUnidades 60
Set-Fila -Fila 234 -Fecha 44511 -Variedad "Sin especificar"    -Calidad "Primera" -Volumen 80  -PrecioMin 10000 -PrecioMax 10000 -PrecioProm 10000 -Unidad "`$/caja 60 unidades" -Origen "Limache"                      -PrecioKg 167 -KgUnidades 60
Set-Fila -Fila 235 -Fecha 44511 -Variedad "Sin especificar"    -Calidad "Primera" -Volumen 240 -PrecioMin 7000  -PrecioMax 8000  -PrecioProm 7500  -Unidad "`$/caja 60 unidades" -Origen "Región de Arica y Parinacota" -PrecioKg 125 -KgUnidades 60
Set-Fila -Fila 236 -Fecha 44511 -Variedad "Sin especificar"    -Calidad "Primera" -Volumen 400 -PrecioMin 9000  -PrecioMax 10000 -PrecioProm 9500  -Unidad "`$/caja 60 unidades" -Origen "Región de O'Higgins"          -PrecioKg 158 -KgUnidades 60
Set-Fila -Fila 237 -Fecha 44511 -Variedad "Sin especificar"    -Calidad "Primera" -Volumen 100 -PrecioMin 9000  -PrecioMax 10000 -PrecioProm 9500  -Unidad "`$/caja 60 unidades" -Origen "Región del Maule"             -PrecioKg 158 -KgUnidades 60

# Make sure the date cells keep the workbook's date display style.
$ws.Range("D233:D237").Style = $ws.Range("D238").Style
